# "contingencies with rene fine"
#
# Two new rows (line7, line8) are inserted into the data table right after
# the existing line1..line6 rows and before the extr1..extr8 rows. That
# pushes extr1..extr8 down by two rows (their index column A is
# renumbered accordingly), and a handful of their from_bus/to_bus/
# in_service values change along the way.
#
# We avoid Range.Insert() here because it fabricates a brand-new cell
# style (duplicating the existing bold/border style minus the border) for
# the freshly inserted cells instead of reusing the original style index.
# Row-by-row Copy(Destination) keeps the original style table untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 8..15 ("extr1".."extr8") down to rows 10..17,
# working bottom-to-top so we never clobber a row before it's copied.
# (Copy(Destination) is used instead of Copy()+PasteSpecial() because
# pasting past the sheet's original used-range boundary with PasteSpecial
# silently drops the source cell style here.)
for ($r = 15; $r -ge 8; $r--) {
    $ws.Range("A" + $r + ":E" + $r).Copy($ws.Range("A" + ($r + 2) + ":E" + ($r + 2)))
}

# New row 8: line7
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# New row 9: line8
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Updated index column + values for the shifted extr1..extr8 rows
# (A renumbers 6..13 -> 8..15; a few from_bus/to_bus/in_service values change)
$extr = @(
    @(10, 8,  5,  12, $true),
    @(11, 9,  5,  9,  $true),
    @(12, 10, 10, 11, $true),
    @(13, 11, 7,  8,  $false),
    @(14, 12, 9,  11, $false),
    @(15, 13, 7,  11, $false),
    @(16, 14, 5,  7,  $true),
    @(17, 15, 8,  5,  $false)
)
foreach ($row in $extr) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
